$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "92.909.09"
$ws.Range("E2").Value = "  -1.68%  "

# Row 3
$ws.Range("D3").Value = "3.366.24"
$ws.Range("E3").Value = "  -1.90%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").Value = "233.17"
$ws.Range("E5").Value = "  -1.83%  "

# Row 6
$ws.Range("D6").Value = "616.62"
$ws.Range("E6").Value = "  -3.97%  "

# Row 7
$ws.Range("D7").Value = "1.37"
$ws.Range("E7").Value = "  -5.10%  "

# Row 8
$ws.Range("D8").Value = "0.389"
$ws.Range("E8").Value = "  -3.72%  "

# Row 9
$ws.Range("E9").Value = "  +0.04%  "

# Row 10
$ws.Range("D10").Value = "0.948"
$ws.Range("E10").Value = "  -2.38%  "

# Row 11
$ws.Range("D11").Value = "3.367.25"
$ws.Range("E11").Value = "  -1.80%  "

# Row 12
$ws.Range("D12").Value = "42.80"
$ws.Range("E12").Value = "  +2.17%  "

# Row 13
$ws.Range("D13").Value = "0.197"
$ws.Range("E13").Value = "  -0.83%  "

# Row 14
$ws.Range("D14").Value = "6.23"
$ws.Range("E14").Value = "  +0.57%  "

# Row 15
$ws.Range("D15").Value = "92.783.30"
$ws.Range("E15").Value = "  -1.61%  "

# Row 16
$ws.Range("D16").Value = "4.006.26"
$ws.Range("E16").Value = "  -1.68%  "

# Row 17
$ws.Range("E17").Value = "  -2.21%  "

# Row 18
$ws.Range("D18").Value = "8.11"
$ws.Range("E18").Value = "  -1.86%  "

# Row 19
$ws.Range("D19").Value = "3.365.71"
$ws.Range("E19").Value = "  -2.03%  "

# Row 20
$ws.Range("D20").Value = "17.41"
$ws.Range("E20").Value = "  -0.74%  "

# Row 21
$ws.Range("D21").Value = "11.26"
$ws.Range("E21").Value = "  -1.09%  "

# Row 22
$ws.Range("E22").Value = "  +4.29%  "

# Row 23
$ws.Range("D23").Value = "495.60"
$ws.Range("E23").Value = "  -0.69%  "

# Row 24
$ws.Range("D24").Value = "0.429"
$ws.Range("E24").Value = "  -15.40%  "

# Row 25
$ws.Range("E25").Value = "  +2.05%  "

# Row 26
$ws.Range("E26").Value = "  -4.84%  "

# Row 27
$ws.Range("D27").Value = "90.74"
$ws.Range("E27").Value = "  -3.78%  "

# Row 28
$ws.Range("D28").Value = "11.99"
$ws.Range("E28").Value = "  +0.85%  "

# Row 29
$ws.Range("D29").Value = "3.541.30"
$ws.Range("E29").Value = "  -2.14%  "

# Row 30
$ws.Range("E30").Value = "  +0.05%  "

# Row 31
$ws.Range("D31").Value = "11.11"
$ws.Range("E31").Value = "  -4.75%  "

# Row 32
$ws.Range("E32").Value = "  -1.50%  "

# Row 33
$ws.Range("D33").Value = "2.68"
$ws.Range("E33").Value = "  -3.20%  "

# Row 34
$ws.Range("D34").Value = "1.01"
$ws.Range("E34").Value = "  +1.04%  "

# Row 35
$ws.Range("D35").Value = "0.172"
$ws.Range("E35").Value = "  -2.35%  "

# Row 36
$ws.Range("D36").Value = "28.57"
$ws.Range("E36").Value = "  -4.25%  "

# Row 37
$ws.Range("E37").Value = "  -4.78%  "

# Row 38
$ws.Range("D38").Value = "557.64"
$ws.Range("E38").Value = "  +1.59%  "

# Row 39
$ws.Range("D39").Value = "7.49"
$ws.Range("E39").Value = "  -1.79%  "

# Row 40
$ws.Range("E40").Value = "  +0.03%  "

# Row 41
$ws.Range("E41").Value = "  -1.15%  "

# Row 42
$ws.Range("E42").Value = "  -4.31%  "

# Row 43
$ws.Range("D43").Value = "0.881"
$ws.Range("E43").Value = "  -1.86%  "

# Row 44
$ws.Range("D44").Value = "23.68"
$ws.Range("E44").Value = "  -1.48%  "

# Row 45
$ws.Range("D45").Value = "1.71"
$ws.Range("E45").Value = "  +0.32%  "

# Row 46
$ws.Range("D46").Value = "3.59"
$ws.Range("E46").Value = "  +0.68%  "

# Row 47
$ws.Range("D47").Value = "0.0407"
$ws.Range("E47").Value = "  -0.21%  "

# Row 48
$ws.Range("D48").Value = "5.41"
$ws.Range("E48").Value = "  -2.54%  "

# Row 49
$ws.Range("D49").Value = "2.11"
$ws.Range("E49").Value = "  -2.61%  "

# Row 50
$ws.Range("D50").Value = "52.62"
$ws.Range("E50").Value = "  -2.95%  "

# Row 51 - Cosmos -> Fantom
$ws.Range("B51").Value = "Fantom"
$ws.Range("C51").Value = "https://coinranking.com/coin/uIEWfMFnQo9K_+fantom-ftm"
$ws.Range("D51").Value = "1.12"
$ws.Range("E51").Value = "  +17.44%  "
